$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5035300850868225
$ws.Range("B1").Value = 0.7527438998222351
$ws.Range("C1").Value = 4.007334232330322
$ws.Range("D1").Value = 2.066445827484131
$ws.Range("E1").Value = 1.072391033172607
